$wb = $excel.ActiveWorkbook

# "ca_hoc" is the second sheet (tab index 1)
$caHoc = $wb.Worksheets.Item(2)

# Update the class-time ranges (buoi hoc) in rows 5-7
$caHoc.Range("A5").Value = "14:00:00"
$caHoc.Range("B5").Value = "15:30:00"
$caHoc.Range("A6").Value = "15:30:00"
$caHoc.Range("B6").Value = "17:00:00"
$caHoc.Range("A7").Value = "20:00:00"
$caHoc.Range("B7").Value = "21:30:00"

# Switch the active sheet to "ca_hoc" and select B9 there
$caHoc.Activate()
$caHoc.Range("B9").Select()
